# Adding test cases in IAM module
#
# The "FindProfileWithInterestTest" row (row 18) on the "Test Cases" sheet
# currently shows a Results value of "PASS". Update it to "SKIP" to reflect
# the latest test run status.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("D18").Value = "SKIP"
